$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "Result" in column D, matching the style of the other header cells
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value = "Result"

# Fill column D (rows 2-13) with alternating Pass/Fail values, matching the
# style used by the other data cells in the row
for ($i = 2; $i -le 13; $i++) {
    if ($i % 2 -eq 0) {
        $value = "Pass"
    } else {
        $value = "Fail"
    }
    $ws.Cells.Item($i, 3).Copy()
    $ws.Cells.Item($i, 4).PasteSpecial(-4122)
    $ws.Cells.Item($i, 4).Value = $value
}

$excel.CutCopyMode = $false

# Update the active selection to match the edited workbook
$ws.Range("D6").Select()
